# BSTATS-105: All teams in season upload for one league (just teams not players)
# Add an "Abbreviation" header above the existing team-abbreviation column (A),
# widen column A to fit it, and leave the selection on the first data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Abbreviation"

$ws.Columns.Item(1).ColumnWidth = 21.5

$ws.Range("A3").Select() | Out-Null
